$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Comp"
$ws.Range("C2").Value = "Itga5"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.4181183333333333
$ws.Range("H2").Value = 1.254355
$ws.Range("I2").Value = 0.01571144052599341
$ws.Range("J2").Value = 0.01571144052599341
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 28.72417333333333
$ws.Range("N2").Value = 86.17251999999999
$ws.Range("O2").Value = 0.4233259107972328
$ws.Range("P2").Value = 0.4233259107972328
$ws.Range("Q2").Value = 12.01010348051111
$ws.Range("R2").Value = 108.0909313246
$ws.Range("S2").Value = 0.006651059870602713
$ws.Range("T2").Value = 0.006651059870602713

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Comp"
$ws.Range("C3").Value = "Itga5"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.4181183333333333
$ws.Range("H3").Value = 1.254355
$ws.Range("I3").Value = 0.01571144052599341
$ws.Range("J3").Value = 0.01571144052599341
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 30.56986233333333
$ws.Range("N3").Value = 91.709587
$ws.Range("O3").Value = 0.4505269713084062
$ws.Range("P3").Value = 0.4505269713084062
$ws.Range("Q3").Value = 12.78181988904278
$ws.Range("R3").Value = 115.036379001385
$ws.Range("S3").Value = 0.007078427715067963
$ws.Range("T3").Value = 0.007078427715067962

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Comp"
$ws.Range("C4").Value = "Itga5"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.4181183333333333
$ws.Range("H4").Value = 1.254355
$ws.Range("I4").Value = 0.01571144052599341
$ws.Range("J4").Value = 0.01571144052599341
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 8.559531999999999
$ws.Range("N4").Value = 25.678596
$ws.Range("O4").Value = 0.126147117894361
$ws.Range("P4").Value = 0.126147117894361
$ws.Range("Q4").Value = 3.578897253953333
$ws.Range("R4").Value = 32.21007528558
$ws.Range("S4").Value = 0.001981952940322732
$ws.Range("T4").Value = 0.001981952940322731

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Comp"
$ws.Range("C5").Value = "Itga5"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 24.89087966666667
$ws.Range("H5").Value = 74.672639
$ws.Range("I5").Value = 0.9353131502385497
$ws.Range("J5").Value = 0.9353131502385496
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 28.72417333333333
$ws.Range("N5").Value = 86.17251999999999
$ws.Range("O5").Value = 0.4233259107972328
$ws.Range("P5").Value = 0.4233259107972328
$ws.Range("Q5").Value = 714.9699419644754
$ws.Range("R5").Value = 6434.72947768028
$ws.Range("S5").Value = 0.395942291205363
$ws.Range("T5").Value = 0.3959422912053631

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Comp"
$ws.Range("C6").Value = "Itga5"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 24.89087966666667
$ws.Range("H6").Value = 74.672639
$ws.Range("I6").Value = 0.9353131502385497
$ws.Range("J6").Value = 0.9353131502385496
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 30.56986233333333
$ws.Range("N6").Value = 91.709587
$ws.Range("O6").Value = 0.4505269713084062
$ws.Range("P6").Value = 0.4505269713084062
$ws.Range("Q6").Value = 760.9107647655659
$ws.Range("R6").Value = 6848.196882890094
$ws.Range("S6").Value = 0.4213838008018981
$ws.Range("T6").Value = 0.4213838008018981

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Comp"
$ws.Range("C7").Value = "Itga5"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 24.89087966666667
$ws.Range("H7").Value = 74.672639
$ws.Range("I7").Value = 0.9353131502385497
$ws.Range("J7").Value = 0.9353131502385496
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 8.559531999999999
$ws.Range("N7").Value = 25.678596
$ws.Range("O7").Value = 0.126147117894361
$ws.Range("P7").Value = 0.126147117894361
$ws.Range("Q7").Value = 213.0542810149826
$ws.Range("R7").Value = 1917.488529134844
$ws.Range("S7").Value = 0.1179870582312885
$ws.Range("T7").Value = 0.1179870582312885

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Comp"
$ws.Range("C8").Value = "Itga5"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1.303350666666667
$ws.Range("H8").Value = 3.910052
$ws.Range("I8").Value = 0.04897540923545694
$ws.Range("J8").Value = 0.04897540923545693
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 28.72417333333333
$ws.Range("N8").Value = 86.17251999999999
$ws.Range("O8").Value = 0.4233259107972328
$ws.Range("P8").Value = 0.4233259107972328
$ws.Range("Q8").Value = 37.43767046344889
$ws.Range("R8").Value = 336.93903417104
$ws.Range("S8").Value = 0.02073255972126701
$ws.Range("T8").Value = 0.02073255972126701

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Comp"
$ws.Range("C9").Value = "Itga5"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1.303350666666667
$ws.Range("H9").Value = 3.910052
$ws.Range("I9").Value = 0.04897540923545694
$ws.Range("J9").Value = 0.04897540923545693
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 30.56986233333333
$ws.Range("N9").Value = 91.709587
$ws.Range("O9").Value = 0.4505269713084062
$ws.Range("P9").Value = 0.4505269713084062
$ws.Range("Q9").Value = 39.84325045205823
$ws.Range("R9").Value = 358.589254068524
$ws.Range("S9").Value = 0.02206474279144016
$ws.Range("T9").Value = 0.02206474279144016

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Comp"
$ws.Range("C10").Value = "Itga5"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1.303350666666667
$ws.Range("H10").Value = 3.910052
$ws.Range("I10").Value = 0.04897540923545694
$ws.Range("J10").Value = 0.04897540923545693
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 8.559531999999999
$ws.Range("N10").Value = 25.678596
$ws.Range("O10").Value = 0.126147117894361
$ws.Range("P10").Value = 0.126147117894361
$ws.Range("Q10").Value = 11.15607173855467
$ws.Range("R10").Value = 100.404645646992
$ws.Range("S10").Value = 0.006178106722749761
$ws.Range("T10").Value = 0.00617810672274976

